# daily auto push: 2026-01-19 13:53 UTC
# A new reading for 2026/01/19 (月, hour 19) was captured and appended to
# the already-existing "today" block, ahead of the pre-populated future
# schedule rows. Insert a new row at 671 and shift the rest down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(671).Insert()

# Force plain-text so the date-like string isn't auto-converted to a real
# Excel date (matches how every other date cell in column A is stored).
$ws.Range("A671").NumberFormat = "@"
$ws.Range("A671").Value = "2026/01/19"
$ws.Range("A671").ClearFormats()

$ws.Range("B671").Value = "月"
$ws.Range("C671").Value = 19
$ws.Range("D671").Value = 201
